$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 263, shifting existing rows 263:317 down to 264:318
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row 263 with the new weekly record for
# Vega Modelo de Temuco / Granada / Wonderfull / Primera / Provincia de Limari
$ws.Cells.Item(263, 1).Value = 10
$ws.Cells.Item(263, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(263, 3).Value = "La Araucanía"
$ws.Cells.Item(263, 4).Value = 45244
$ws.Cells.Item(263, 5).Value = 9
$ws.Cells.Item(263, 6).Value = "Fruta"
$ws.Cells.Item(263, 7).Value = 100104
$ws.Cells.Item(263, 8).Value = "Frutos de pepita"
$ws.Cells.Item(263, 9).Value = 100104001
$ws.Cells.Item(263, 10).Value = "Granada"
$ws.Cells.Item(263, 11).Value = "Wonderfull"
$ws.Cells.Item(263, 12).Value = "Primera"
$ws.Cells.Item(263, 13).Value = 50
$ws.Cells.Item(263, 14).Value = 18000
$ws.Cells.Item(263, 15).Value = 18000
$ws.Cells.Item(263, 16).Value = 18000
$ws.Cells.Item(263, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(263, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(263, 19).Value = 1800
$ws.Cells.Item(263, 20).Value = 10

# Ensure the date cell keeps the same number format as the other date cells in column D
$ws.Cells.Item(263, 4).NumberFormat = $ws.Cells.Item(264, 4).NumberFormat
